# Applies the "Updated capital structure database" edit described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: refresh existing first data row values ---
$ws.Cells.Item(2, 1).Value = "Kazakhstan"
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "3"
$ws.Cells.Item(2, 3).Value = "Bank (Money Center)"
$ws.Cells.Item(2, 4).Value = 0.185
$ws.Cells.Item(2, 5).Value = 0.234
$ws.Cells.Item(2, 7).Value = 0.0
$ws.Cells.Item(2, 8).Value = 0.0
$ws.Cells.Item(2, 9).Value = 0.0
$ws.Cells.Item(2, 10).Value = 0.0
$ws.Cells.Item(2, 11).Value = 893.8
$ws.Cells.Item(2, 12).Value = 0.5204681767891457
$ws.Cells.Item(2, 13).Value = 560.51
$ws.Cells.Item(2, 14).Value = 0.1371144108221825
$ws.Cells.Item(2, 15).Value = 0.6271089729245917
$ws.Cells.Item(2, 16).Value = 510.9
$ws.Cells.Item(2, 17).Value = 0.1249785953668143
$ws.Cells.Item(2, 18).Value = 0.5716043857686284
$ws.Cells.Item(2, 19).Value = 49.60999999999999
$ws.Cells.Item(2, 20).Value = 0.08850867959536848
$ws.Cells.Item(2, 21).Value = 6149.4
$ws.Cells.Item(2, 22).Value = 1.504293157856112
$ws.Cells.Item(2, 23).Value = 0.2197253433208489
$ws.Cells.Item(2, 24).Value = 0.08477942022756023
$ws.Cells.Item(2, 25).Value = 0.1349459230932887
$ws.Cells.Item(2, 26).Value = 0.4441626991270379
$ws.Cells.Item(2, 27).Value = 0.0
$ws.Cells.Item(2, 28).Value = 0.06660264414975592
$ws.Cells.Item(2, 29).Value = -0.06660264414975592
$ws.Cells.Item(2, 30).Value = 4148.9
$ws.Cells.Item(2, 31).Value = 0.0
$ws.Cells.Item(2, 32).Value = 4148.9
$ws.Cells.Item(2, 33).Value = -2000.5
$ws.Cells.Item(2, 34).Value = 0.5037028943278943
$ws.Cells.Item(2, 35).Value = 0.5060004390565163
$ws.Cells.Item(2, 36).Value = -0.9583692631982372
$ws.Cells.Item(2, 37).Value = -0.9758536585365853
$ws.Cells.Item(2, 38).Value = 0.0
$ws.Cells.Item(2, 39).Value = 0.0

# Columns AN2 (debt_ebitda) and AP2 (net_debt_ebitda) are no longer populated
$ws.Cells.Item(2, 40).ClearContents()
$ws.Cells.Item(2, 42).ClearContents()

# --- Row 3: Joint Stock Company Halyk Savings Bank of Kazakhstan (LSE:HSBK) ---
# Clear the old debt_ebitda / net_debt_ebitda values that used to live on this row
$ws.Cells.Item(3, 40).ClearContents()
$ws.Cells.Item(3, 42).ClearContents()
$ws.Cells.Item(3, 1).Value = "Kazakhstan"
$ws.Cells.Item(3, 2).Value = "Joint Stock Company Halyk Savings Bank of Kazakhstan (LSE:HSBK)"
$ws.Cells.Item(3, 3).Value = "Bank (Money Center)"
$ws.Cells.Item(3, 4).Value = 0.201
$ws.Cells.Item(3, 5).Value = 0.234
$ws.Cells.Item(3, 7).Value = 0.0
$ws.Cells.Item(3, 8).Value = 0.0
$ws.Cells.Item(3, 9).Value = 0.0
$ws.Cells.Item(3, 10).Value = 0.0
$ws.Cells.Item(3, 11).Value = 757.5
$ws.Cells.Item(3, 12).Value = 0.5483567395395975
$ws.Cells.Item(3, 13).Value = 509.4
$ws.Cells.Item(3, 14).Value = 0.1494849898758694
$ws.Cells.Item(3, 15).Value = 0.6724752475247524
$ws.Cells.Item(3, 16).Value = 461.7
$ws.Cells.Item(3, 17).Value = 0.1354872788097544
$ws.Cells.Item(3, 18).Value = 0.6095049504950495
$ws.Cells.Item(3, 19).Value = 47.69999999999999
$ws.Cells.Item(3, 20).Value = 0.09363957597173143
$ws.Cells.Item(3, 21).Value = 4945.7
$ws.Cells.Item(3, 22).Value = 1.451330809636999
$ws.Cells.Item(3, 23).Value = 0.2395408405274642
$ws.Cells.Item(3, 24).Value = 0.06209217875287636
$ws.Cells.Item(3, 25).Value = 0.1774486617745878
$ws.Cells.Item(3, 26).Value = 0.654521333301746
$ws.Cells.Item(3, 27).Value = 0.0
$ws.Cells.Item(3, 28).Value = 0.05478503337207604
$ws.Cells.Item(3, 29).Value = -0.05478503337207604
$ws.Cells.Item(3, 30).Value = 2587.7
$ws.Cells.Item(3, 31).Value = 0.0
$ws.Cells.Item(3, 32).Value = 2587.7
$ws.Cells.Item(3, 33).Value = -2358.0
$ws.Cells.Item(3, 34).Value = 0.4316142375821463
$ws.Cells.Item(3, 35).Value = 0.4501365526118948
$ws.Cells.Item(3, 36).Value = -2.246356101743356
$ws.Cells.Item(3, 37).Value = -2.936488169364882
$ws.Cells.Item(3, 38).Value = 0.0
$ws.Cells.Item(3, 39).Value = 0.0

# --- Row 4 (new): ForteBank Joint Stock Company (KAS:ASBN) ---
$ws.Cells.Item(4, 1).Value = "Kazakhstan"
$ws.Cells.Item(4, 2).Value = "ForteBank Joint Stock Company (KAS:ASBN)"
$ws.Cells.Item(4, 3).Value = "Bank (Money Center)"
$ws.Cells.Item(4, 4).Value = 0.185
$ws.Cells.Item(4, 5).Value = -0.243
$ws.Cells.Item(4, 7).Value = 0.0
$ws.Cells.Item(4, 8).Value = 0.0
$ws.Cells.Item(4, 9).Value = 0.0
$ws.Cells.Item(4, 10).Value = 0.0
$ws.Cells.Item(4, 11).Value = 123.2
$ws.Cells.Item(4, 12).Value = 0.5317220543806647
$ws.Cells.Item(4, 13).Value = 51.11
$ws.Cells.Item(4, 14).Value = 0.08639283299526707
$ws.Cells.Item(4, 15).Value = 0.4148538961038961
$ws.Cells.Item(4, 16).Value = 49.2
$ws.Cells.Item(4, 17).Value = 0.08316430020283976
$ws.Cells.Item(4, 18).Value = 0.3993506493506493
$ws.Cells.Item(4, 19).Value = 1.909999999999997
$ws.Cells.Item(4, 20).Value = 0.03737037761690465
$ws.Cells.Item(4, 21).Value = 629.0
$ws.Cells.Item(4, 22).Value = 1.063218390804598
$ws.Cells.Item(4, 23).Value = 0.2197253433208489
$ws.Cells.Item(4, 24).Value = 0.08477942022756023
$ws.Cells.Item(4, 25).Value = 0.1349459230932887
$ws.Cells.Item(4, 26).Value = 0.2050442477876106
$ws.Cells.Item(4, 27).Value = 0.0
$ws.Cells.Item(4, 28).Value = 0.06660264414975592
$ws.Cells.Item(4, 29).Value = -0.06660264414975592
$ws.Cells.Item(4, 30).Value = 960.1
$ws.Cells.Item(4, 31).Value = 0.0
$ws.Cells.Item(4, 32).Value = 960.1
$ws.Cells.Item(4, 33).Value = 331.1
$ws.Cells.Item(4, 34).Value = 0.618740735967004
$ws.Cells.Item(4, 35).Value = 0.6199793361746094
$ws.Cells.Item(4, 36).Value = 0.3588381922618403
$ws.Cells.Item(4, 37).Value = 0.3600478468899522
$ws.Cells.Item(4, 38).Value = 0.0
$ws.Cells.Item(4, 39).Value = 0.0

# --- Row 5 (new): Joint Stock Company Bank CenterCredit (KAS:CCBN) ---
$ws.Cells.Item(5, 1).Value = "Kazakhstan"
$ws.Cells.Item(5, 2).Value = "Joint Stock Company Bank CenterCredit (KAS:CCBN)"
$ws.Cells.Item(5, 3).Value = "Bank (Money Center)"
$ws.Cells.Item(5, 4).Value = 0.091
$ws.Cells.Item(5, 5).Value = 1.45
$ws.Cells.Item(5, 7).Value = 0.0
$ws.Cells.Item(5, 8).Value = 0.0
$ws.Cells.Item(5, 9).Value = 0.0
$ws.Cells.Item(5, 10).Value = 0.0
$ws.Cells.Item(5, 11).Value = 13.1
$ws.Cells.Item(5, 12).Value = 0.1257197696737044
$ws.Cells.Item(5, 13).Value = 0.0
$ws.Cells.Item(5, 14).Value = 0.0
$ws.Cells.Item(5, 15).Value = 0.0
$ws.Cells.Item(5, 16).Value = 0.0
$ws.Cells.Item(5, 17).Value = 0.0
$ws.Cells.Item(5, 18).Value = 0.0
$ws.Cells.Item(5, 19).Value = 0.0
$ws.Cells.Item(5, 21).Value = 574.7
$ws.Cells.Item(5, 22).Value = 6.48645598194131
$ws.Cells.Item(5, 23).Value = 0.04235370190753314
$ws.Cells.Item(5, 24).Value = 0.2203885934186463
$ws.Cells.Item(5, 25).Value = -0.1780348915111132
$ws.Cells.Item(5, 26).Value = 0.1664999536612413
$ws.Cells.Item(5, 27).Value = 0.0
$ws.Cells.Item(5, 28).Value = 0.08224434324894245
$ws.Cells.Item(5, 29).Value = -0.08224434324894245
$ws.Cells.Item(5, 30).Value = 601.1
$ws.Cells.Item(5, 31).Value = 0.0
$ws.Cells.Item(5, 32).Value = 601.1
$ws.Cells.Item(5, 33).Value = 26.39999999999998
$ws.Cells.Item(5, 34).Value = 0.8715383500072496
$ws.Cells.Item(5, 35).Value = 0.6663341093005211
$ws.Cells.Item(5, 36).Value = 0.2295652173913042
$ws.Cells.Item(5, 37).Value = 0.08063530849114227
$ws.Cells.Item(5, 38).Value = 0.0
$ws.Cells.Item(5, 39).Value = 0.0

